# Fruta / hortaliza, semanal
# The weekly refresh reshuffles which date-row holds which set of
# market-report figures (columns D, M-T). Column A-L (market/product
# metadata) stay the same for every row in this sheet, only the
# date/volume/price/unit/origin columns are redistributed across rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-row figures being redistributed.
$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the "before" values for every affected row so that writes
# (which happen row by row) never read already-overwritten data.
$snapshot = @{}
for ($r = 2; $r -le 26; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Destination row -> source row (which row's original data now lives here).
$rowMap = @{
    2  = 26
    3  = 25
    4  = 10
    5  = 16
    6  = 15
    7  = 9
    8  = 20
    9  = 21
    10 = 5
    11 = 13
    12 = 12
    13 = 7
    14 = 2
    15 = 3
    16 = 19
    17 = 18
    18 = 17
    19 = 24
    20 = 11
    21 = 6
    22 = 14
    23 = 22
    24 = 8
    25 = 4
    26 = 23
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
